$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove column M entirely; the old column N (with its data) shifts left
# to become the new column M.
$ws.Range("M:M").EntireColumn.Delete()

# Move/confirm selection on the (now last) column M, row 1, matching the
# post-edit workbook state.
$ws.Range("M1").Select()
